# Update InsideBet Data: Automatizado
# Applies the latest Bundesliga standings update to the three affected rows:
#  - Row 10 (Hamburger SV): refreshed match-week stats
#  - Row 14 / Row 15 (Gladbach / Mainz 05): teams swap table positions,
#    Mainz 05 also receives refreshed stats for the newly played match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Hamburger SV ---
$ws.Range("C10").Value = 22
$ws.Range("E10").Value = 8
$ws.Range("G10").Value = 25
$ws.Range("H10").Value = 32
$ws.Range("J10").Value = 26
$ws.Range("K10").Value = 1.18
$ws.Range("L10").Value = "D D W W D"

# --- Row 14: becomes Mainz 05 (with updated stats) ---
$ws.Range("B14").Value = "Mainz 05"
$ws.Range("C14").Value = 23
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 7
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 26
$ws.Range("H14").Value = 38
$ws.Range("I14").Value = -12
$ws.Range("J14").Value = 22
$ws.Range("K14").Value = 0.96
$ws.Range("L14").Value = "W W W L D"
$ws.Range("M14").Value = 30939
$ws.Range("N14").Value = "Nadiem Amiri - 10"
$ws.Range("O14").Value = "Daniel Batz"

# --- Row 15: becomes Gladbach (unchanged stats from previous position) ---
$ws.Range("B15").Value = "Gladbach"
$ws.Range("C15").Value = 22
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 10
$ws.Range("G15").Value = 25
$ws.Range("H15").Value = 37
$ws.Range("I15").Value = -12
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = "D L D D L"
$ws.Range("M15").Value = 50747
$ws.Range("N15").Value = "Haris Tabakovic - 10"
$ws.Range("O15").Value = "Moritz Nicolas"
